$d = $word.ActiveDocument
$matchCase = $true
$matchWholeWord = $true
$matchWildcards = $false
$matchSoundsLike = $false
$matchAllWordForms = $false
$forward = $true
$wrap = 1
$format = $false
$replaceOne = 1

$d.Content.Find.Execute("2025-12-11 Thursday", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "2025-12-12 Friday", $replaceOne) | Out-Null
$d.Content.Find.Execute("90-5=85", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "93-56=37", $replaceOne) | Out-Null
$d.Content.Find.Execute("81-8=73", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "7+36=43", $replaceOne) | Out-Null
$d.Content.Find.Execute("68+9=77", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "98-79=19", $replaceOne) | Out-Null
$d.Content.Find.Execute("73-64=9", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "54+27=81", $replaceOne) | Out-Null
$d.Content.Find.Execute("29+15=44", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "64+27=91", $replaceOne) | Out-Null
$d.Content.Find.Execute("58+39=97", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "58+15=73", $replaceOne) | Out-Null
$d.Content.Find.Execute("47+15=62", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "90-51=39", $replaceOne) | Out-Null
$d.Content.Find.Execute("80-45=35", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "7+47=54", $replaceOne) | Out-Null
$d.Content.Find.Execute("8+63=71", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "9+8=17", $replaceOne) | Out-Null
$d.Content.Find.Execute("22-8=14", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "92-73=19", $replaceOne) | Out-Null
$d.Content.Find.Execute("39+2=41", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "13+29=42", $replaceOne) | Out-Null
$d.Content.Find.Execute("81-12=69", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "9+43=52", $replaceOne) | Out-Null
$d.Content.Find.Execute("85-17=68", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "69+16=85", $replaceOne) | Out-Null
$d.Content.Find.Execute("61-49=12", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "36+26=62", $replaceOne) | Out-Null
$d.Content.Find.Execute("19+3=22", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "8+44=52", $replaceOne) | Out-Null
$d.Content.Find.Execute("95-56=39", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "43-26=17", $replaceOne) | Out-Null
$d.Content.Find.Execute("67+14=81", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "14+17=31", $replaceOne) | Out-Null
$d.Content.Find.Execute("82-78=4", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "66+29=95", $replaceOne) | Out-Null
$d.Content.Find.Execute("90-71=19", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "54+17=71", $replaceOne) | Out-Null
$d.Content.Find.Execute("9+16=25", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "93-49=44", $replaceOne) | Out-Null
$d.Content.Find.Execute("17+64=81", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "68+5=73", $replaceOne) | Out-Null
$d.Content.Find.Execute("60-14=46", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "16+66=82", $replaceOne) | Out-Null
$d.Content.Find.Execute("92-79=13", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "44+19=63", $replaceOne) | Out-Null
$d.Content.Find.Execute("48+4=52", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "37+8=45", $replaceOne) | Out-Null
$d.Content.Find.Execute("86-9=77", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "24+58=82", $replaceOne) | Out-Null
$d.Content.Find.Execute("73-46=27", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "20-1=19", $replaceOne) | Out-Null
$d.Content.Find.Execute("93-34=59", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "56+7=63", $replaceOne) | Out-Null
$d.Content.Find.Execute("92-27=65", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "90-42=48", $replaceOne) | Out-Null
$d.Content.Find.Execute("42-37=5", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "39+33=72", $replaceOne) | Out-Null
$d.Content.Find.Execute("19+47=66", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "67-49=18", $replaceOne) | Out-Null
$d.Content.Find.Execute("34+49=83", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "28+9=37", $replaceOne) | Out-Null
$d.Content.Find.Execute("92-77=15", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "76+6=82", $replaceOne) | Out-Null
$d.Content.Find.Execute("18+43=61", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "82-18=64", $replaceOne) | Out-Null
$d.Content.Find.Execute("22-9=13", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "33+18=51", $replaceOne) | Out-Null
$d.Content.Find.Execute("68+19=87", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "44+29=73", $replaceOne) | Out-Null
$d.Content.Find.Execute("44+8=52", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "95-68=27", $replaceOne) | Out-Null
$d.Content.Find.Execute("48+38=86", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "37+39=76", $replaceOne) | Out-Null
$d.Content.Find.Execute("34+28=62", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "93-55=38", $replaceOne) | Out-Null
$d.Content.Find.Execute("97-68=29", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "84-48=36", $replaceOne) | Out-Null
$d.Content.Find.Execute("40-34=6", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "97-8=89", $replaceOne) | Out-Null
$d.Content.Find.Execute("80-67=13", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "74-7=67", $replaceOne) | Out-Null
$d.Content.Find.Execute("81-18=63", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "7+68=75", $replaceOne) | Out-Null
$d.Content.Find.Execute("24+9=33", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "64+27=91", $replaceOne) | Out-Null
$d.Content.Find.Execute("72-16=56", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "47+26=73", $replaceOne) | Out-Null
$d.Content.Find.Execute("64-55=9", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "14-8=6", $replaceOne) | Out-Null
$d.Content.Find.Execute("96-89=7", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "13+19=32", $replaceOne) | Out-Null
$d.Content.Find.Execute("4+7=11", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "29+55=84", $replaceOne) | Out-Null
$d.Content.Find.Execute("15-7=8", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "8+24=32", $replaceOne) | Out-Null
$d.Content.Find.Execute("82-25=57", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "43-8=35", $replaceOne) | Out-Null
$d.Content.Find.Execute("91-55=36", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "78+15=93", $replaceOne) | Out-Null
$d.Content.Find.Execute("90-7=83", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "15-6=9", $replaceOne) | Out-Null
$d.Content.Find.Execute("18+6=24", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "37+8=45", $replaceOne) | Out-Null
$d.Content.Find.Execute("84-65=19", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "9+84=93", $replaceOne) | Out-Null
$d.Content.Find.Execute("68+14=82", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "7+89=96", $replaceOne) | Out-Null
$d.Content.Find.Execute("84-15=69", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "4+69=73", $replaceOne) | Out-Null
$d.Content.Find.Execute("34+27=61", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "73-38=35", $replaceOne) | Out-Null
$d.Content.Find.Execute("53+8=61", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "19+39=58", $replaceOne) | Out-Null
$d.Content.Find.Execute("36+7=43", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "56-27=29", $replaceOne) | Out-Null
$d.Content.Find.Execute("86-59=27", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "72-54=18", $replaceOne) | Out-Null
$d.Content.Find.Execute("85-8=77", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "84-48=36", $replaceOne) | Out-Null
$d.Content.Find.Execute("90-49=41", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "49+14=63", $replaceOne) | Out-Null
$d.Content.Find.Execute("62-47=15", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "31-24=7", $replaceOne) | Out-Null
$d.Content.Find.Execute("6+39=45", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "84-17=67", $replaceOne) | Out-Null
$d.Content.Find.Execute("36+38=74", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "50-6=44", $replaceOne) | Out-Null
$d.Content.Find.Execute("58+38=96", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "50-16=34", $replaceOne) | Out-Null
$d.Content.Find.Execute("76-27=49", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "27+35=62", $replaceOne) | Out-Null
$d.Content.Find.Execute("62+9=71", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "78+14=92", $replaceOne) | Out-Null
$d.Content.Find.Execute("81-26=55", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "35-29=6", $replaceOne) | Out-Null
$d.Content.Find.Execute("28+36=64", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "30-27=3", $replaceOne) | Out-Null
$d.Content.Find.Execute("95-77=18", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "91-73=18", $replaceOne) | Out-Null
$d.Content.Find.Execute("9+59=68", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "91-45=46", $replaceOne) | Out-Null
$d.Content.Find.Execute("24-15=9", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "41-33=8", $replaceOne) | Out-Null
$d.Content.Find.Execute("28+69=97", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "97-18=79", $replaceOne) | Out-Null
$d.Content.Find.Execute("29+3=32", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "83-69=14", $replaceOne) | Out-Null
$d.Content.Find.Execute("91-86=5", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "18+15=33", $replaceOne) | Out-Null
$d.Content.Find.Execute("62-6=56", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "19+68=87", $replaceOne) | Out-Null
$d.Content.Find.Execute("47+45=92", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "29+42=71", $replaceOne) | Out-Null
$d.Content.Find.Execute("76-48=28", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "73-8=65", $replaceOne) | Out-Null
$d.Content.Find.Execute("35+58=93", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "44+38=82", $replaceOne) | Out-Null
$d.Content.Find.Execute("43-5=38", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "37+25=62", $replaceOne) | Out-Null
$d.Content.Find.Execute("23-4=19", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "15+27=42", $replaceOne) | Out-Null
$d.Content.Find.Execute("49+27=76", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "96-19=77", $replaceOne) | Out-Null
$d.Content.Find.Execute("5+57=62", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "70-6=64", $replaceOne) | Out-Null
$d.Content.Find.Execute("42+49=91", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "64+19=83", $replaceOne) | Out-Null
$d.Content.Find.Execute("65-6=59", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "46+17=63", $replaceOne) | Out-Null
$d.Content.Find.Execute("25+17=42", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "43+48=91", $replaceOne) | Out-Null
$d.Content.Find.Execute("96-69=27", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "61-5=56", $replaceOne) | Out-Null
$d.Content.Find.Execute("9+22=31", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "79+5=84", $replaceOne) | Out-Null
$d.Content.Find.Execute("19+3=22", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "38+7=45", $replaceOne) | Out-Null
$d.Content.Find.Execute("94-15=79", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "43-16=27", $replaceOne) | Out-Null
$d.Content.Find.Execute("34-16=18", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "77+4=81", $replaceOne) | Out-Null
$d.Content.Find.Execute("23+59=82", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "42+19=61", $replaceOne) | Out-Null
$d.Content.Find.Execute("7+8=15", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "93-7=86", $replaceOne) | Out-Null
$d.Content.Find.Execute("66+26=92", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "62-35=27", $replaceOne) | Out-Null
$d.Content.Find.Execute("67+26=93", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "88+3=91", $replaceOne) | Out-Null
$d.Content.Find.Execute("40-3=37", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "75+17=92", $replaceOne) | Out-Null
$d.Content.Find.Execute("9+85=94", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "45-16=29", $replaceOne) | Out-Null
$d.Content.Find.Execute("8+27=35", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "83-37=46", $replaceOne) | Out-Null
$d.Content.Find.Execute("37+34=71", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "17+68=85", $replaceOne) | Out-Null
$d.Content.Find.Execute("10-8=2", $matchCase, $matchWholeWord, $matchWildcards, $matchSoundsLike, $matchAllWordForms, $forward, $wrap, $format, "16+16=32", $replaceOne) | Out-Null
